$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.1877
$ws.Range("A9").Value = -20.31009999999997
$ws.Range("A18").Value = -23.00330000000001
$ws.Range("A20").Value = -22.17510000000003
$ws.Range("E21").Value = 13.09289999999999
